# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rule row 11 (B11) on the "Rules" sheet is renamed from "R40" to "1".
# The leading apostrophe forces Excel to store the numeral-looking value
# as literal text (quote-prefixed) instead of coercing it to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
